$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: backpack / tshirt / bikelight / max / planck / 722036
$ws.Range("A2").Value = "backpack"
$ws.Range("B2").Value = "tshirt"
$ws.Range("C2").Value = "bikelight"
$ws.Range("D2").Value = "max"
$ws.Range("E2").Value = "planck"
$ws.Range("F2").Value = 722036
$ws.Range("G2:I2").ClearContents() | Out-Null

# Row 3: (A3 cleared) / tshirt / bikelight / edward / snowden / 753651
$ws.Range("A3").ClearContents() | Out-Null
$ws.Range("B3").Value = "tshirt"
$ws.Range("C3").Value = "bikelight"
$ws.Range("D3").Value = "edward"
$ws.Range("E3").Value = "snowden"
$ws.Range("F3").Value = 753651
$ws.Range("G3:I3").ClearContents() | Out-Null

# Update selection to match G2:I3 with active cell G2
$ws.Range("G2:I3").Select() | Out-Null
